$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number (column B, the hashcode column) to the new hashcode value.
# These correspond to the updated MD5-like hash codes in the "hashcode.csv" sheet.
$updates = @(
    @{Row=43; Value="9908c29b7fab4cef9c80ff46dcf03af6"}
    @{Row=44; Value="5ff9c44479d2673932aeff8998a825e5"}
    @{Row=59; Value="67b7c9f494638cfa165a2cb3d182fc27"}
    @{Row=62; Value="4e05305223e360a8e38718a184190ecd"}
    @{Row=84; Value="9a5db73c9534f604467f508eb15529c7"}
    @{Row=86; Value="6d8341204789516057bbfe961f4c210f"}
    @{Row=88; Value="358cdad6f08cddf8823a1319b7362204"}
    @{Row=89; Value="b2af67e341352a35b5b3c501abeb8995"}
    @{Row=99; Value="62c25067cd76a2f41d1a7e30de2c631d"}
    @{Row=110; Value="6aab2ca26942956d73d9ff8f8c6467f5"}
    @{Row=134; Value="59eecfb103695c4ee96129d3b0bc1abb"}
    @{Row=148; Value="55aa84e32c71003ec1c2da80fa01a948"}
    @{Row=154; Value="f06f4c62a046da5903e0df0fadfdc179"}
    @{Row=161; Value="7b7a2be5778190d36ca610c576f2aba5"}
    @{Row=169; Value="294c74fb396996bcda6a4caf4a9e6251"}
    @{Row=187; Value="ec1698ab8e6353b5c9c7fcc6b4c8e660"}
    @{Row=196; Value="b7d2fb4e6114387d8b68a8ec6efef78b"}
    @{Row=197; Value="c7b13c896d79f110af279a23b4fc0266"}
    @{Row=215; Value="bb9ffa8ad25621f081e6aa7a6eb95859"}
    @{Row=284; Value="c026b5107650650af1812e21e3c1a63c"}
    @{Row=324; Value="b40b1af66372a7b4200f93a3b9a705ab"}
    @{Row=344; Value="852dd907a8027478a4daf60cf9088c2c"}
    @{Row=346; Value="2db18754cfde08b840d87e01bd99babb"}
    @{Row=354; Value="1d538b85277d28472a82b08c1972ef36"}
    @{Row=363; Value="c28618037e25e3704c441d1b6fd012db"}
    @{Row=365; Value="d2366e3876160ff8d83a9104aaef04fe"}
    @{Row=369; Value="ca3bc93378cdc336fb1d23bfd0287161"}
    @{Row=386; Value="88eabee8ff2de015477c019b4a4715d8"}
    @{Row=396; Value="70f3101ce53b40825777f574a3f672eb"}
    @{Row=398; Value="263b79777cba372a1ce8952159ab7e5a"}
    @{Row=401; Value="2e0d91a3db8ecb98665e3c3e1cba4772"}
    @{Row=418; Value="3363b4c2e8d7958f3eb77a47037d30e2"}
    @{Row=439; Value="69eaf0d15580db9662dbadd5e7757bf6"}
    @{Row=447; Value="625a89473f71f5ad7c8cf7e7d4c1e5cf"}
    @{Row=469; Value="4e57d31c81b8f0b4e80c3d80a51d4c70"}
    @{Row=493; Value="3dcb045817c3098469dbe3b3069eb83d"}
    @{Row=497; Value="1651e0a559aaee0fdaf68eecd726d95c"}
    @{Row=519; Value="198651d299d2370b52a79c886970ff2e"}
    @{Row=538; Value="935ef81b9fb805a4be34a828255ed1c9"}
    @{Row=549; Value="7da4e1704fa9188165cebf1be3b85e08"}
    @{Row=552; Value="b28a09c762509387e269714acae243b2"}
    @{Row=569; Value="65d97a694b9a1bafb5141fa8bcf07fdf"}
    @{Row=571; Value="e6ece49d19fe9fdfa605f0790e42ec84"}
    @{Row=615; Value="a2749027397b979e75718696e5075f86"}
    @{Row=668; Value="61f5ab1a750febd2448d6d4ca0bbfe4e"}
    @{Row=669; Value="87d77c62805a11f381beb9c5e15f72ae"}
    @{Row=699; Value="31cfbe580fe6dfa303d681e53f909b80"}
    @{Row=722; Value="7392593711ab7bb3250d7e56e6dece16"}
    @{Row=725; Value="b63bbf7509084a2b7b03ebc9663565f9"}
    @{Row=733; Value="6bf45b8ee91ad771ecde83c2df66c768"}
    @{Row=735; Value="ea3f030a8a490b970f7da01d292fecd4"}
    @{Row=748; Value="f0e5c24140c9cb41727b2a54e3a260c8"}
    @{Row=761; Value="ed8d0f2e534d224e43bcd71563f0f88b"}
    @{Row=763; Value="a7c1dadd7ec65360b7f5725a4cbe5b52"}
    @{Row=781; Value="ba34bad9b6e7c800a1ade4e878ae97ca"}
    @{Row=785; Value="6ef8f9ed2d01b92bcc647bec0e721f5d"}
    @{Row=818; Value="dfb41bc87eb854e2b99e152486144459"}
    @{Row=828; Value="99ff8c6d85ab8542a3e97d9d6608161d"}
    @{Row=832; Value="65c35ad342e492fe7c438999db77b425"}
    @{Row=833; Value="cf1807c266eaf5e8f2d88606fba8fae9"}
    @{Row=837; Value="b4e3ba5466c57f698549f42028666b21"}
    @{Row=840; Value="d8ef158d83f87a66c3dcaa3efb83d98e"}
    @{Row=844; Value="012a4a70a89fbaee5703d8a2e3ea5c6e"}
    @{Row=918; Value="b2f44d3255fdd42d13a8b4353660a499"}
    @{Row=931; Value="b213c4e2488ad106b7235cb5839dce7a"}
    @{Row=940; Value="95060d6e31bc91f2529e80a514b7f8a5"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.Value
}
